$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.345.48'
$ws.Range('E2').Value = '  -2.23%  '
$ws.Range('D3').Value = '3.492.09'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '612.47'
$ws.Range('E5').Value = '  +4.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.21'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').Value = '  +1.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -2.68%  '
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.26'
$ws.Range('E11').Value = '  -2.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000307'
$ws.Range('E12').Value = '  -3.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.59'
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('D14').Value = '4.051.61'
$ws.Range('E14').Value = '  -2.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '605.59'
$ws.Range('E15').Value = '  +4.30%  '
$ws.Range('D16').Value = '69.431.63'
$ws.Range('E16').Value = '  -2.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '12.66'
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.79'
$ws.Range('E18').Value = '  -2.62%  '
$ws.Range('D19').Value = '3.477.92'
$ws.Range('E19').Value = '  -2.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.988'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.31'
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '105.22'
$ws.Range('E23').Value = '  +10.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.66'
$ws.Range('E24').Value = '  +1.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.02'
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('E26').Value = '  +2.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.95'
$ws.Range('E27').Value = '  -3.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('E28').Value = '  +9.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.72'
$ws.Range('E29').Value = '  +3.76%  '
$ws.Range('E30').Value = '  -3.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.48'
$ws.Range('E31').Value = '  +1.47%  '
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('E33').Value = '  +15.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.31'
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('E35').Value = '  -6.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '522.85'
$ws.Range('E37').Value = '  -5.10%  '
$ws.Range('E38').Value = '  -5.09%  '
$ws.Range('D39').Value = '3.577.12'
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.58'
$ws.Range('E40').Value = '  +4.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.79'
$ws.Range('E41').Value = '  -2.60%  '
$ws.Range('E42').Value = '  -3.64%  '
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0461'
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('E45').Value = '  +1.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.144'
$ws.Range('E46').Value = '  +5.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.32'
$ws.Range('E47').Value = '  -5.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.85'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '131.01'
$ws.Range('E50').Value = '  -3.78%  '
$ws.Range('E51').Value = '  -9.62%  '
